$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of column L header cells into new M/N header cells
# and set their values to "Pattern" and "Pattern Type" for each of the three
# property-type blocks (rows 4, 12, 20).
$headerRows = @(4, 12, 20)
foreach ($r in $headerRows) {
    $lCell = $ws.Cells.Item($r, 12)  # L
    $mCell = $ws.Cells.Item($r, 13)  # M
    $nCell = $ws.Cells.Item($r, 14)  # N

    $lCell.Copy()
    $mCell.PasteSpecial(-4122)  # xlPasteFormats
    $nCell.PasteSpecial(-4122)  # xlPasteFormats

    $mCell.Value = "Pattern"
    $nCell.Value = "Pattern Type"
}

$ws.Range("M20:N20").Select()
